# Generate Report for Handback
# Adds a new handback-status row for 15400665-184c-48eb-8bf3-2f4a8b20e85b.md
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileId   = "15400665-184c-48eb-8bf3-2f4a8b20e85b"
$fileName = "$fileId.md"
$pathName = "e2e\$fileId.md"
$status   = "Handed back: in sync with en-US"
$ext      = ".md"

$zhXlf = "$fileId.1e512fd1b2e4f149a71b39182ae048e392d24200.zh-cn.xlf"
$deXlf = "$fileId.1e512fd1b2e4f149a71b39182ae048e392d24200.de-de.xlf"

$dtHandoff   = "2016-08-14 17:00:28"
$dtZhHoDate  = "2016-08-14 17:00:20"
$dtZhHbDate  = "2016-08-14 17:00:48"
$dtDeHbDate  = "2016-08-14 17:00:59"

$srcRepoUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/$fileName"
$zhRepoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/$fileName"
$deRepoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/$fileName"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> table "Overview" (A1:G3 -> A1:G4)
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A4").Value = $fileName
$wsOv.Range("B4").Value = $pathName
$wsOv.Range("C4").Value = $ext
$wsOv.Range("E4").Value = $status
$wsOv.Range("F4").Value = $status
$wsOv.Range("G4").Value = $dtHandoff

$wsOv.Range("B4").Style = "Hyperlink"
$wsOv.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), $srcRepoUrl, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> table "zh_cn" (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $dtZhHoDate
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $dtZhHbDate
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Range("A4").Style = "Hyperlink"
$wsZh.Range("I4").Style = "Hyperlink"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcRepoUrl, "", "", $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhRepoUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" -> table "de_de" (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $dtHandoff
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $dtDeHbDate
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Range("A4").Style = "Hyperlink"
$wsDe.Range("I4").Style = "Hyperlink"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcRepoUrl, "", "", $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deRepoUrl, "", "", $fileName) | Out-Null

Write-Output "Handback row added for $fileName"
